# expansão das análises automáticas
# Adds three new summary columns (apoio_medio, contribuicoes, media_contribuicoes)
# to the sheet, and rescales the existing "particip" / "taxa_sucesso" columns
# (E:F) from fractions to percentage-like numbers (value * 100), keeping their
# existing percentage display format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (row 1), matching the header style used by A1:K1 ---
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)   # xlPasteFormats: copy formatting only
$excel.CutCopyMode = $false

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- Rescale existing columns E (particip) and F (taxa_sucesso): x100 ---
$ws.Range("E2").Value = 82.02247191011236
$ws.Range("F2").Value = 61.7351598173516

$ws.Range("E3").Value = 17.97752808988764
$ws.Range("F3").Value = 64.16666666666667

$ws.Range("E4").Value = 75.54495912806539
$ws.Range("F4").Value = 93.05680793507665

$ws.Range("E5").Value = 24.45504087193461
$ws.Range("F5").Value = 97.7715877437326

$ws.Range("E6").Value = 90.64327485380117
$ws.Range("F6").Value = 20.64516129032258

$ws.Range("E7").Value = 9.35672514619883
$ws.Range("F7").Value = 37.5

# --- New data columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes) ---
$ws.Range("L2").Value = 90.1471541030692
$ws.Range("M2").Value = 211660
$ws.Range("N2").Value = 313.1065088757396

$ws.Range("L3").Value = 96.01937043276359
$ws.Range("M3").Value = 51893
$ws.Range("N3").Value = 336.9675324675325

$ws.Range("L4").Value = 88.30982863725519
$ws.Range("M4").Value = 147585
$ws.Range("N4").Value = 143.0087209302326

$ws.Range("L5").Value = 95.05584770392593
$ws.Range("M5").Value = 56061
$ws.Range("N5").Value = 159.7179487179487

$ws.Range("L6").Value = 17.70131305034959
$ws.Range("M6").Value = 1877
$ws.Range("N6").Value = 14.6640625

$ws.Range("L7").Value = 30.09544761111365
$ws.Range("M7").Value = 331
$ws.Range("N7").Value = 13.79166666666667
